$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Row 57 (T15_4, first data row) ---
$ws.Range("D57").Value = 720
$ws.Range("E57").Value = 4400
$ws.Range("F57").Value = 24164
$ws.Range("G57").Value = 0.04
$ws.Range("H57").Value = 16
$ws.Range("I57").Value = 0.36

# --- Row 58 ---
$ws.Range("D58").Value = 1260
$ws.Range("E58").Value = 7665
$ws.Range("F58").Value = 24194
$ws.Range("G58").Value = 0.06
$ws.Range("H58").Value = 21
$ws.Range("I58").Value = 0.48

# --- Row 59 ---
$ws.Range("D59").Value = 2760
$ws.Range("E59").Value = 16713
$ws.Range("F59").Value = 24220
$ws.Range("G59").Value = 0.15
$ws.Range("H59").Value = 33
$ws.Range("I59").Value = 1.06

# --- Row 60 ---
$ws.Range("D60").Value = 4680
$ws.Range("E60").Value = 28273
$ws.Range("F60").Value = 24234
$ws.Range("G60").Value = 0.21
$ws.Range("H60").Value = 50
$ws.Range("I60").Value = 1.36

# --- Row 61 ---
$ws.Range("D61").Value = 8360
$ws.Range("E61").Value = 50408
$ws.Range("F61").Value = 24240
$ws.Range("G61").Value = 0.4
$ws.Range("H61").Value = 83
$ws.Range("I61").Value = 3.38

# Clear the old "T15_8" label that used to live at B66 (it moves up to B64)
$ws.Range("B66").ClearContents()

# --- Row 64 (new label row: T15_8) ---
$ws.Range("B64").Value = "T15_8"
$ws.Range("D64").Value = 720
$ws.Range("E64").Value = 4481
$ws.Range("F64").Value = 24248
$ws.Range("G64").Value = 0.07
$ws.Range("H64").Value = 20
$ws.Range("I64").Value = 0.45

# --- Row 65 ---
$ws.Range("D65").Value = 1260
$ws.Range("E65").Value = 7771
$ws.Range("F65").Value = 24250
$ws.Range("G65").Value = 0.11
$ws.Range("H65").Value = 28
$ws.Range("I65").Value = 0.84

# --- Row 66 ---
$ws.Range("D66").Value = 2760
$ws.Range("E66").Value = 16867
$ws.Range("F66").Value = 24250
$ws.Range("G66").Value = 0.33
$ws.Range("H66").Value = 49
$ws.Range("I66").Value = 1.91

# --- Row 67 ---
$ws.Range("D67").Value = 4680
$ws.Range("E67").Value = 28467
$ws.Range("F67").Value = 24252
$ws.Range("G67").Value = 0.41
$ws.Range("H67").Value = 78
$ws.Range("I67").Value = 2.52

# --- Row 68 (overwrite old values, used to hold the 912/5667/... row) ---
$ws.Range("D68").Value = 8360
$ws.Range("E68").Value = 50657
$ws.Range("F68").Value = 24252
$ws.Range("G68").Value = 0.67
$ws.Range("H68").Value = 133
$ws.Range("I68").Value = 4.03

# --- Row 69 no longer holds data; clear it entirely ---
$ws.Range("D69:I69").ClearContents()

# Update the sheet view to match: scrolled down so row 57 is the top row,
# with the final selection sitting on I70.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 57
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I70").Select()
